$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the internal role cell identifiers (testRun -> testRunItem)
$ws.Range("A1").Value = "div_testRunItem_internalRoleCellName"
$ws.Range("B1").Value = "div_testRunItem_internalRoleCellName_1"

# Widen columns A (34 -> 38) and B (36 -> 40).
# COM's ColumnWidth setter/getter round-trips through a +5/6 character-width
# pad relative to the raw OOXML <col width> value, so back it out here to
# land on the exact stored widths.
$pad = 5 / 6
$ws.Columns.Item(1).ColumnWidth = 38 - $pad
$ws.Columns.Item(2).ColumnWidth = 40 - $pad
